$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7..91 down to 8..92
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new record
$ws.Cells.Item(7, 1).Value  = 11
$ws.Cells.Item(7, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value  = "Bíobío"
$ws.Cells.Item(7, 4).Value  = 44537
$ws.Cells.Item(7, 5).Value  = 8
$ws.Cells.Item(7, 6).Value  = 100112043
$ws.Cells.Item(7, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(7, 8).Value  = "Sin especificar"
$ws.Cells.Item(7, 9).Value  = "Primera"
$ws.Cells.Item(7, 10).Value = 220
$ws.Cells.Item(7, 11).Value = 9000
$ws.Cells.Item(7, 12).Value = 10000
$ws.Cells.Item(7, 13).Value = 9545
$ws.Cells.Item(7, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(7, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 16).Value = 159
$ws.Cells.Item(7, 17).Value = 60
$ws.Cells.Item(7, 18).Value = "Hortaliza"
